$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws.Range("D2").Value = 5046.52
$ws.Range("E2").Value = -5046.52

$ws.Range("D4").Value = 6378.320000000001
$ws.Range("E4").Value = 7345.02
$ws.Range("F4").Value = 0.4647789823760106
